$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column B (shifts existing B:U to D:W)
$ws.Range("B:C").Insert()

# Set header values for the new columns
$ws.Range("B1").Value = "rest_begin"
$ws.Range("C1").Value = "rest_end"

# Set data values for the new columns
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 162

# Update selection to match target state
$ws.Range("C3").Select()
